# Tweak slides and remove old file
#
# 1) Slide 10 speaker notes: extend the existing sentence about Visual
#    Studio/VS Code with a new sentence about JetBrains Rider's community
#    license.
# 2) Slide 9 speaker notes: reword the "Explore the integrated development
#    environments..." sentence.

$p = $ppt.ActivePresentation

# --- Slide 10 notes: Visual Studio / JetBrains Rider blurb -----------------
$s10 = $p.Slides.Item(10)
$notes10 = $s10.NotesPage
$notesShape10 = $notes10.Shapes.Item(2)
$notesShape10.TextFrame.TextRange.Text = "Visual Studio on Windows and macOS are primary IDEs, complemented by VS Code and command-line tools that support .NET MAUI development workflows. JetBrains Rider is also available on both platforms for folks who either have a subscription or use it through the community license."

# --- Slide 9 notes: IDEs/tools intro sentence -------------------------------
$s9 = $p.Slides.Item(9)
$notes9 = $s9.NotesPage
$notesShape9 = $notes9.Shapes.Item(2)
$notesShape9.TextFrame.TextRange.Text = "Next we’ll discuss the IDEs and tools available on Windows and macOS to build, debug, and deploy .NET MAUI applications."
